# Generate Report for Handoff
# A new handoff was produced for file "7472a77c-e5fd-479e-a3a9-956a736b8f16" (row 4 in the
# per-locale detail sheets). Its "Latest Handoff Datetime" is refreshed for both locales, and
# the roll-up "Latest Handoff Date" on the Overview sheet follows suit.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-20 16:46:18"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-20 16:46:26"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D4").Value = "2016-03-20 16:46:26"
